$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 108 (shifts existing rows 108-175 down to 109-176,
# and expands the used range to A1:R176).
$ws.Rows(108).Insert()

# Populate the newly inserted row 108 with the new record's data.
$ws.Range("A108").Value = 10
$ws.Range("B108").Value = "Vega Modelo de Temuco"
$ws.Range("C108").Value = "La Araucanía"
$ws.Range("D108").Value = 44596
$ws.Range("E108").Value = 9
$ws.Range("F108").Value = 100112005
$ws.Range("G108").Value = "Puerro"
$ws.Range("H108").Value = "Azul de Maquehue"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 10
$ws.Range("K108").Value = 13000
$ws.Range("L108").Value = 13000
$ws.Range("M108").Value = 13000
$ws.Range("N108").Value = "$/docena de paquetes"
$ws.Range("O108").Value = "Provincia de Cautín"
$ws.Range("P108").Value = 1083
$ws.Range("Q108").Value = 12
$ws.Range("R108").Value = "Hortaliza"
